$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (week number + date range) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Crime-stat table updates (rows 14-30) ---
$ws.Range("C14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 11.764705882352
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = -17.355371900826
$ws.Range("L16").Value = -28.571428571428
$ws.Range("M16").Value = -32.432432432432
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -35.714285714285
$ws.Range("I17").Value = 150
$ws.Range("J17").Value = 181
$ws.Range("K17").Value = -17.127071823204
$ws.Range("L17").Value = -18.918918918918
$ws.Range("M17").Value = 31.578947368421
$ws.Range("C18").Value = 2
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = 20.967741935483
$ws.Range("L18").Value = -14.772727272727
$ws.Range("M18").Value = 27.118644067796
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 15.625
$ws.Range("I19").Value = 213
$ws.Range("J19").Value = 196
$ws.Range("K19").Value = 8.673469387755
$ws.Range("L19").Value = 13.297872340425
$ws.Range("M19").Value = 33.125
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 97
$ws.Range("J20").Value = 94
$ws.Range("K20").Value = 3.191489361702
$ws.Range("L20").Value = 155.263157894737
$ws.Range("M20").Value = 136.585365853659
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = -2.061855670103
$ws.Range("I21").Value = 645
$ws.Range("J21").Value = 665
$ws.Range("K21").Value = -3.007518796992
$ws.Range("L21").Value = -0.769230769230
$ws.Range("M21").Value = 20.560747663551
$ws.Range("D22").Value = 1
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = 0
$ws.Range("H22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 9
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -55
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 12.5
$ws.Range("C23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("H23").Value = "'***.*"
$ws.Range("A23").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 53.333333333333
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = -1.111111111111
$ws.Range("I24").Value = 599
$ws.Range("J24").Value = 906
$ws.Range("K24").Value = -33.885209713024
$ws.Range("L24").Value = -34.392113910186
$ws.Range("M24").Value = 72.126436781609
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 9.090909090909
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -13.157894736842
$ws.Range("I25").Value = 252
$ws.Range("J25").Value = 263
$ws.Range("K25").Value = -4.182509505703
$ws.Range("L25").Value = 1.204819277108
$ws.Range("M25").Value = -17.105263157894
$ws.Range("D26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'***.*"
$ws.Range("A26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("G26").Value = 2
$ws.Range("L26").Value = -46.666666666666
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 200
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 25
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = -28.571428571428
$ws.Range("L27").Value = -10.714285714285
$ws.Range("C28").Value = "'0"
$ws.Range("A28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").Value = "'0"
$ws.Range("A29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("I30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("K30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("I30").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("K30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = -33.333333333333
